$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "id" / "product_id" / "image_id" / "cart_id" / "user_id" labels to
# --- their PK/FK annotated versions (model de cart y user cambiados)

# PC (products) mini-table header row 3
$ws.Range("E3").Value = "id(PK)"
$ws.Range("O3").Value = "id(PK)"
$ws.Range("P3").Value = "product_id(FK)"
$ws.Range("Q3").Value = "image_id(FK)"
$ws.Range("S3").Value = "id(PK)"

# CART_PRODUCT / CARTS header row 15
$ws.Range("E15").Value = "id(PK)"
$ws.Range("J15").Value = "id(PK)"
$ws.Range("K15").Value = "product_id(FK)"
$ws.Range("L15").Value = "cart_id(FK)"
$ws.Range("N15").Value = "id(PK)"
$ws.Range("P15").Value = "user_id(FK)"

# USERS header row 27
$ws.Range("E27").Value = "id(PK)"

# --- Selection moved from S20 to N15 ---
$ws.Range("N15").Select()

# --- New column widths for columns P (16) and Q (17) ---
# (engine quantizes ColumnWidth to 1/6-character steps; these inputs land on
#  the closest achievable stored widths to the target 14.140625 / 12.28515625)
$ws.Columns.Item(16).ColumnWidth = 13.35
$ws.Columns.Item(17).ColumnWidth = 11.5
